# Update "想去人数" (F column) counts on both the "展览" sheet and the
# aggregated "全部类型" sheet, which both contain the same rows of data.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Mapping of row number -> new value for column F
$updates = @{
    3  = 581
    10 = 5096
    12 = 15
    13 = 33
    14 = 5
    15 = 46
    16 = 180
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
